$wb = $excel.ActiveWorkbook

# --- "March" sheet: insert two new expense rows above the "Monthly total" row ---
$march = $wb.Worksheets.Item("March")

# Current layout (before edit):
#   Row 1: headers (Category, Name, Date, Price, Account)
#   Rows 2-9: expense entries
#   Row 10: Monthly total row (A10 = "Monthly total: ", B10 = 18360.0)
#
# New layout (after edit):
#   Rows 2-9: unchanged expense entries
#   Row 10: new expense entry (Transportation / asdf / 2023-03-22 / 1000.0 / Checkings)
#   Row 11: new expense entry (Entertainment / douchebag / 2023-03-22 / 3000.0 / Checkings)
#   Row 12: Monthly total row (A12 = "Monthly total: ", B12 = 22360.0)

# Insert two new rows right before the old total row (row 10); this pushes the
# existing total row (and its values) down to row 12 automatically.
$march.Rows.Item(10).Insert()
$march.Rows.Item(10).Insert()

# All data in this workbook is stored as text (even things that look like
# numbers/dates), so force text formatting before assigning "numeric-looking"
# or "date-looking" strings to avoid Excel auto-converting them, then clear
# the format again so no stray style survives on the cell.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# New row 10
$march.Range("A10").Value = "Transportation"
$march.Range("B10").Value = "asdf"
Set-TextValue $march.Range("C10") "2023-03-22"
Set-TextValue $march.Range("D10") "1000.0"
$march.Range("E10").Value = "Checkings"

# New row 11
$march.Range("A11").Value = "Entertainment"
$march.Range("B11").Value = "douchebag"
Set-TextValue $march.Range("C11") "2023-03-22"
Set-TextValue $march.Range("D11") "3000.0"
$march.Range("E11").Value = "Checkings"

# Row 12 is the (shifted down) monthly total row; update its value to reflect
# the two newly added expenses (18360 + 1000 + 3000 = 22360).
$march.Range("A12").Value = "Monthly total: "
$march.Range("B12").Value = 22360.0

# --- "April" sheet: no content change, the total row label is untouched ---
$april = $wb.Worksheets.Item("April")
$april.Range("A3").Value = "Monthly total: "
